$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (req 1.1): Admin (not Manager) registers employees ---
$ws.Range("B2").Value = "System/company Admin"
$ws.Range("B2").WrapText = $true
$ws.Range("D2").Value = "Admin can register an employee adding them to the system"

# --- Row 7 (req 1.6): Admin (not Manager) removes employees ---
$ws.Range("B7").Value = "System/company Admin"
$ws.Range("B7").WrapText = $true
$ws.Range("D7").Value = "Admin can remove/make inactive an employee from the system"

# --- Row 14 (req 3.3): clarify wording + fix related-reqs cell to a real line break ---
$ws.Range("D14").Value = "Manager can override/correct employee hours for past and future (misreport, sick, coming vacation, etc)"
$ws.Range("F14").Value = "'2.2" + [char]10 + "2.4"
$ws.Range("F14").WrapText = $true

# Stray formatted-but-empty cell next to it (matches the author's worksheet exactly)
$ws.Range("G14").Value = "'"
$ws.Range("G14").WrapText = $true
$ws.Range("G14").ClearContents()

# --- Row 15 (req 3.4): typo fix confimred -> confirmed ---
$ws.Range("D15").Value = "System notifies if Manager is overriding existing value and checks to ensure override is confirmed"

# --- Formatting cleanup: wrap the whole Description column ---
$ws.Range("D2:D17").WrapText = $true

# --- Misc view state ---
$ws.PageSetup.Orientation = 1
[void]$ws.Range("H13").Select()

Write-Output "done"
